$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column H (shifts H: onward right by 2)
$ws.Range("H1:I1").EntireColumn.Insert()

# Set header values for the two new columns
$ws.Range("H1").Value = "FUENTE"
$ws.Range("I1").Value = "SUBFUENTE"

# Adjust column widths (values chosen so the host's internal pixel
# rounding lands as close as possible to the target stored widths)
$ws.Range("E1").ColumnWidth = 19.666666666666664
$ws.Range("F1").ColumnWidth = 13.666666666666668
$ws.Range("I1").ColumnWidth = 15.0

# Reflect the author's final selection/scroll position on the sheet
[void]$ws.Range("I1").Select()
$excel.ActiveWindow.ScrollColumn = 3
